# ez_hubs.xlsx - "ez lc upload fix"
# Adds four new active/ocean hub rows (Dalian, Hong Kong, Bahrain, Shanghai)
# to the hubs sheet, including the Shanghai photo hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 7: Dalian -----------------------------------------------------
$ws.Range("A7").Value = "active"
$ws.Range("B7").Value = "ocean"
$ws.Range("C7").Value = "Dalian"
$ws.Range("E7").Value = 38.9310971
$ws.Range("F7").Value = 121.6590698
$ws.Range("G7").Value = "China"
$ws.Range("H7").Value = "Yuejin Rd, Zhongshan Qu, Dalian Shi, Liaoning Sheng, China"

# ---- Row 8: Hong Kong ----------------------------------------------------
$ws.Range("A8").Value = "active"
$ws.Range("B8").Value = "ocean"
$ws.Range("C8").Value = "Hong Kong"
$ws.Range("E8").Value = 22.3081225
$ws.Range("F8").Value = 114.220595
$ws.Range("G8").Value = "China"
$ws.Range("H8").Value = "Hong Kong, Kwun Tong, 偉業街223-231號,宏利金融中心B座7字樓,701A及708B室"

# ---- Row 9: Bahrain -----------------------------------------------------
$ws.Range("A9").Value = "active"
$ws.Range("B9").Value = "ocean"
$ws.Range("C9").Value = "Bahrain"
$ws.Range("E9").Value = 26.1480232
$ws.Range("F9").Value = 50.6433176
$ws.Range("G9").Value = "Bahrain"
$ws.Range("H9").Value = "Sitra, Bahrain"

# ---- Row 10: Shanghai (with photo + hyperlink) --------------------------
$ws.Range("A10").Value = "active"
$ws.Range("B10").Value = "ocean"
$ws.Range("C10").Value = "Shanghai"
$ws.Range("D10").Font.Name = "Arial"
$ws.Range("E10").Value = 30.626539
$ws.Range("F10").Value = 122.064958
$ws.Range("G10").Value = "China"
$ws.Range("H10").Value = "Zhoushan, Zhejiang, China"
$ws.Range("I10").Value = "https://assets.itsmycargo.com/assets/cityimages/Shanghai_sm.jpg"
$ws.Hyperlinks.Add($ws.Range("I10"), "https://assets.itsmycargo.com/assets/cityimages/Shanghai_sm.jpg")
$ws.Range("I10").Font.Name = "Arial"
$ws.Range("I10").Font.Underline = $true
$ws.Range("I10").Font.Color = 13391121

# Trailing formatted-but-empty cells that mirror the header row's styled
# run (J:Y), left over from the row-10 paste in the original edit.
for ($col = 10; $col -le 25; $col++) {
    $ws.Cells.Item(10, $col).Font.Name = "Arial"
}
